# "Duplicate Daily Data Fix iOS"
# The daily iOS-ratings pull re-ran (2023-06-27 22:30:24 -> 2023-06-28 10:50:42),
# so every data row's Detail Date / Date stamp advances by a day, and the
# review-count / rank / rating columns pick up the freshly scraped figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDetailDate = "2023-06-28 10:50:42"
$newDate       = "June 28, 2023"

# Per-row overrides for iOS Total Reviews (F), iOS Rank (G) and, for row 34,
# iOS App Rating (E) -- only cells whose numbers actually moved are listed.
$rowUpdates = @{
    2  = @{ F = 61798 }
    5  = @{ F = 399366;  G = 101 }
    6  = @{ F = 4755313 }
    7  = @{ F = 43177;   G = 96 }
    8  = @{ F = 322339 }
    9  = @{ F = 2412836; G = 4 }
    10 = @{ F = 42793;   G = 92 }
    11 = @{ F = 24997;   G = 71 }
    12 = @{ F = 2011630 }
    13 = @{ F = 1047294 }
    14 = @{ F = 258155;  G = 113 }
    15 = @{ F = 116111;  G = 189 }
    16 = @{ F = 77293;   G = 152 }
    18 = @{ F = 411965;  G = 53 }
    19 = @{ F = 3991280; G = 14 }
    20 = @{ G = 158 }
    22 = @{ G = 59 }
    23 = @{ G = 192 }
    24 = @{ F = 877815 }
    25 = @{ F = 303 }
    26 = @{ F = 1625 }
    32 = @{ F = 1377 }
    33 = @{ F = 28397 }
    34 = @{ E = 3.1; F = 9 }
    35 = @{ F = 13446 }
}

for ($row = 2; $row -le 38; $row++) {
    $ws.Range("B$row").Value = $newDetailDate
    $ws.Range("C$row").Value = $newDate

    if ($rowUpdates.ContainsKey($row)) {
        $cols = $rowUpdates[$row]
        foreach ($col in $cols.Keys) {
            $ws.Range("$col$row").Value = $cols[$col]
        }
    }
}
